$wb = $excel.ActiveWorkbook
$wsQSBA = $wb.Worksheets.Item("QSBA")
$wsQSDA = $wb.Worksheets.Item("QSDA")

# Set column G (rows 2-69) to 0 on the QSBA sheet.
# The QSDA sheet's G/I columns use formulas (VLOOKUP against QSBA!E:G and
# IF(I>0,I,H)) so they will recalculate automatically once these values change.
for ($r = 2; $r -le 69; $r++) {
    $wsQSBA.Cells.Item($r, 7).Value = 0
}

# Update the QSBA sheet's remembered scroll/selection state: row 46 becomes
# the top-left visible row and G2 is the selected cell.
$wsQSBA.Activate()
$excel.ActiveWindow.ScrollRow = 46
$wsQSBA.Range("G2").Select() | Out-Null

# QSDA remains the active/visible sheet (tabSelected), restore it.
$wsQSDA.Activate()

$excel.Calculate()
